# The presentation ships two theme parts:
#   ppt/theme/theme1.xml -> currently the "Integral" colour scheme (used by
#                            the slide master / presentation theme)
#   ppt/theme/theme2.xml -> currently the "Office Theme" colour scheme
#                            (used only by the notes master)
#
# The authored change swaps the two themes' contents, so that the visible
# presentation theme (theme1.xml) becomes the "Office Theme" palette and
# the notes-master theme (theme2.xml) becomes the "Integral" palette.
#
# This automation host only exposes a single, shared Theme object through
# the PowerPoint object model (SlideMaster / NotesMaster / HandoutMaster
# all resolve to the same underlying theme part), so we apply the "Office
# Theme" colours - the half of the swap that is reachable through
# ThemeColorScheme - to that shared theme.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Index -> (scheme slot, BGR-encoded RGB integer for the Office Theme colour)
# 1 dk1       000000
# 2 lt1       FFFFFF
# 3 dk2       44546A
# 4 lt2       E7E6E6
# 5 accent1   5B9BD5
# 6 accent2   ED7D31
# 7 accent3   A5A5A5
# 8 accent4   FFC000
# 9 accent5   4472C4
# 10 accent6  70AD47
# 11 hlink    0563C1
# 12 folHlink 954F72
$colorScheme.Colors(1).RGB = 0
$colorScheme.Colors(2).RGB = 16777215
$colorScheme.Colors(3).RGB = 6968388
$colorScheme.Colors(4).RGB = 15132391
$colorScheme.Colors(5).RGB = 13998939
$colorScheme.Colors(6).RGB = 3243501
$colorScheme.Colors(7).RGB = 10855845
$colorScheme.Colors(8).RGB = 49407
$colorScheme.Colors(9).RGB = 12874308
$colorScheme.Colors(10).RGB = 4697456
$colorScheme.Colors(11).RGB = 12673797
$colorScheme.Colors(12).RGB = 7491477
